$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Insert three new paragraphs - the "Development Sprint 3" heading,
#    the "Test Case 1" body copy, and a (relocated) empty paragraph
#    holding the _GoBack bookmark - right before the "Links:" heading
#    that follows the Sprint 2 / Test Case 2 paragraph.
#
#    Word's Range.InsertXML splices full <w:p> blocks in place of the
#    (possibly zero-length) paragraph the range sits in, so the target
#    paragraph's own XML is appended, verbatim, as the last <w:p> in
#    the inserted fragment - that paragraph effectively "replaces
#    itself" and the new paragraphs land immediately before it.
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Links:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the 'Links:' paragraph"
}
$linksPara = $rng.Paragraphs(1)

$insertXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Development Sprint 3—Cloud Apps Test Bench</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr><w:lastRenderedPageBreak/><w:t>Test Case 1:</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> We are going to make sure that the application deploys smoothly onto Carolina Cloud Apps. In order to test this, we will have several users log onto Cloud Apps to ensure that they are able to access the application. </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>We will make sure that these users are a</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">ble to select their availability and that this gets passed to the database. This should then get passed into the scheduling engine and the engine should handle scheduling all of the employees and then return the schedule to the application </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>(the</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> schedules will not be displayed yet).</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00FA646A" w:rsidRDefault="009212E4" w:rsidP="008D419E"><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Links:</w:t></w:r></w:p>
'@

$r = $linksPara.Range
$r.Collapse(1)
$r.InsertXML($insertXml)

# ------------------------------------------------------------------
# 2. The _GoBack bookmark used to live in its own empty paragraph
#    right after the "database/tests" hyperlink further down the
#    document. Since the bookmark now lives in the newly inserted
#    paragraph above, delete that old, now-stale paragraph outright.
# ------------------------------------------------------------------
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldPara = $oldBookmark.Range.Paragraphs(1)
$oldPara.Range.Delete()

Write-Host "done"
